$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 195.30435
$ws.Range("I33").Value = 135.05263
$ws.Range("K33").Value = 135.05263
$ws.Range("M33").Value = 93.94737000000001
$ws.Range("H94").Value = 7037.769
$ws.Range("I94").Value = 3457.5833
$ws.Range("K94").Value = 3457.5833
$ws.Range("M94").Value = -3006.5833
$ws.Range("H103").Value = 907.8889
$ws.Range("I103").Value = 858
$ws.Range("K103").Value = 2574
$ws.Range("M103").Value = -1988
$ws.Range("H106").Value = 1818
$ws.Range("I106").Value = 1745.25
$ws.Range("K106").Value = 1745.25
$ws.Range("M106").Value = -1114.25
$ws.Range("H113").Value = 2910.7778
$ws.Range("J113").Value = 3374.25
$ws.Range("L113").Value = 3374.25
$ws.Range("N113").Value = -9882.25
$ws.Range("H121").Value = 7238
$ws.Range("J121").Value = 7238
$ws.Range("L121").Value = 21714
$ws.Range("N121").Value = -25208
$ws.Range("H132").Value = 139217.17
$ws.Range("I132").Value = 168807.72
$ws.Range("K132").Value = 506423.16
$ws.Range("M132").Value = -503893.16
$ws.Range("H138").Value = 3877.7568
$ws.Range("I138").Value = 1792.4286
$ws.Range("J138").Value = 5147.087
$ws.Range("K138").Value = 5377.2858
$ws.Range("L138").Value = 15441.261
$ws.Range("M138").Value = -237.2857999999997
$ws.Range("N138").Value = -25721.261
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""
$ws.Range("H141").Value = 1609.25
$ws.Range("I141").Value = 1609.25
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4827.75
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 352.25
$ws.Range("N141").Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 16534.467
$ws.Range("I2").Value = 6885.7144
$ws.Range("K2").Value = 6885.7144
$ws.Range("M2").Value = -6772.7144
$ws.Range("H61").Value = 3032178.2
$ws.Range("I61").Value = 3704995.8
$ws.Range("K61").Value = 3704995.8
$ws.Range("M61").Value = -3704783.8
$ws.Range("H74").Value = 2755.75
$ws.Range("I74").Value = 951.0714
$ws.Range("J74").Value = 6966.6665
$ws.Range("K74").Value = 951.0714
$ws.Range("L74").Value = 6966.6665
$ws.Range("M74").Value = -77.07140000000004
$ws.Range("N74").Value = -8714.666499999999
$ws.Range("H77").Value = 2755.75
$ws.Range("I77").Value = 951.0714
$ws.Range("J77").Value = 6966.6665
$ws.Range("K77").Value = 4755.357
$ws.Range("L77").Value = 34833.3325
$ws.Range("M77").Value = -387.357
$ws.Range("N77").Value = -43569.3325
$ws.Range("H116").Value = 16534.467
$ws.Range("I116").Value = 6885.7144
$ws.Range("K116").Value = 6885.7144
$ws.Range("M116").Value = -4591.7144
$ws.Range("H122").Value = 3316.476
$ws.Range("I122").Value = 2547
$ws.Range("K122").Value = 7641
$ws.Range("M122").Value = -5191
$ws.Range("H132").Value = 1178891.4
$ws.Range("I132").Value = 1541013.5
$ws.Range("J132").Value = 1994.75
$ws.Range("K132").Value = 4623040.5
$ws.Range("L132").Value = 5984.25
$ws.Range("M132").Value = -4620510.5
$ws.Range("N132").Value = -11044.25
$ws.Range("H136").Value = 3032178.2
$ws.Range("I136").Value = 3704995.8
$ws.Range("K136").Value = 11114987.4
$ws.Range("M136").Value = -11112437.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 16534.467
$ws.Range("I3").Value = 6885.7144
$ws.Range("K3").Value = 6885.7144
$ws.Range("M3").Value = -6771.7144
$ws.Range("H94").Value = 1169.0968
$ws.Range("I94").Value = 1151.6072
$ws.Range("J94").Value = 1332.3334
$ws.Range("K94").Value = 1151.6072
$ws.Range("L94").Value = 1332.3334
$ws.Range("M94").Value = -700.6071999999999
$ws.Range("N94").Value = -2234.3334
$ws.Range("H134").Value = 1137613
$ws.Range("I134").Value = 1136041.4
$ws.Range("K134").Value = 3408124.2
$ws.Range("M134").Value = -3405589.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 473.6
$ws.Range("I22").Value = 457.33334
$ws.Range("K22").Value = 457.33334
$ws.Range("M22").Value = -107.33334
$ws.Range("H105").Value = 22191.389
$ws.Range("I105").Value = 24228.125
$ws.Range("K105").Value = 24228.125
$ws.Range("M105").Value = -22481.125
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = ""
$ws.Range("H141").Value = 86326
$ws.Range("J141").Value = 86326
$ws.Range("L141").Value = 86326
$ws.Range("N141").Value = -96686

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 766.0833
$ws.Range("I5").Value = 1416.3334
$ws.Range("J5").Value = 115.833336
$ws.Range("K5").Value = 4249.0002
$ws.Range("L5").Value = 347.500008
$ws.Range("M5").Value = -4137.0002
$ws.Range("N5").Value = -571.500008
$ws.Range("H14").Value = 191
$ws.Range("I14").Value = 191
$ws.Range("K14").Value = 573
$ws.Range("M14").Value = -400
$ws.Range("H70").Value = 2625.8
$ws.Range("I70").Value = 1283.5
$ws.Range("K70").Value = 3850.5
$ws.Range("M70").Value = -3535.5
$ws.Range("H73").Value = 2625.8
$ws.Range("I73").Value = 1283.5
$ws.Range("K73").Value = 3850.5
$ws.Range("M73").Value = -2758.5
$ws.Range("H102").Value = 9011.888999999999
$ws.Range("J102").Value = 9011.888999999999
$ws.Range("L102").Value = 27035.667
$ws.Range("N102").Value = -31903.667
$ws.Range("H109").Value = 2530.875
$ws.Range("I109").Value = 1321
$ws.Range("K109").Value = 3963
$ws.Range("M109").Value = -2923
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").Value = ""
$ws.Range("H120").Value = 15000
$ws.Range("I120").Value = 15000
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 45000
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -40162
$ws.Range("N120").Value = ""
$ws.Range("H132").Value = 824.1818
$ws.Range("I132").Value = 792.625
$ws.Range("J132").Value = 908.3333
$ws.Range("K132").Value = 7133.625
$ws.Range("L132").Value = 8174.9997
$ws.Range("M132").Value = -4603.625
$ws.Range("N132").Value = -13234.9997
$ws.Range("H135").Value = 766.0833
$ws.Range("I135").Value = 1416.3334
$ws.Range("J135").Value = 115.833336
$ws.Range("K135").Value = 12747.0006
$ws.Range("L135").Value = 1042.500024
$ws.Range("M135").Value = -10212.0006
$ws.Range("N135").Value = -6112.500024

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 27441.076
$ws.Range("I107").Value = 44091
$ws.Range("J107").Value = 801.2
$ws.Range("K107").Value = 44091
$ws.Range("L107").Value = 801.2
$ws.Range("M107").Value = -42171
$ws.Range("N107").Value = -4641.2
$ws.Range("H132").Value = 59538480
$ws.Range("I132").Value = 92008620
$ws.Range("K132").Value = 276025860
$ws.Range("M132").Value = -276023330

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3876.0476
$ws.Range("I7").Value = 3562.375
$ws.Range("K7").Value = 3562.375
$ws.Range("M7").Value = -3450.375
$ws.Range("H22").Value = 1160
$ws.Range("J22").Value = 1612.5
$ws.Range("L22").Value = 1612.5
$ws.Range("N22").Value = -2202.5
$ws.Range("H27").Value = 1160
$ws.Range("J27").Value = 1612.5
$ws.Range("L27").Value = 1612.5
$ws.Range("N27").Value = -1826.5
$ws.Range("H40").Value = 3738.25
$ws.Range("I40").Value = 3738.25
$ws.Range("K40").Value = 3738.25
$ws.Range("M40").Value = -3602.25
$ws.Range("H122").Value = 3755.4707
$ws.Range("J122").Value = 4758.8
$ws.Range("L122").Value = 14276.4
$ws.Range("N122").Value = -19176.4
$ws.Range("H126").Value = 3876.0476
$ws.Range("I126").Value = 3562.375
$ws.Range("K126").Value = 10687.125
$ws.Range("M126").Value = -8217.125
$ws.Range("H132").Value = 1662418.1
$ws.Range("I132").Value = 2323385.8
$ws.Range("K132").Value = 6970157.399999999
$ws.Range("M132").Value = -6967627.399999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 14998.5
$ws.Range("J8").Value = 14998.5
$ws.Range("L8").Value = 14998.5
$ws.Range("N8").Value = -15278.5
$ws.Range("H33").Value = 18000
$ws.Range("I33").Value = 18000
$ws.Range("K33").Value = 18000
$ws.Range("M33").Value = -17750
$ws.Range("H36").Value = 18000
$ws.Range("I36").Value = 18000
$ws.Range("K36").Value = 18000
$ws.Range("M36").Value = -17750
$ws.Range("H122").Value = 1922.037
$ws.Range("I122").Value = 1619.0952
$ws.Range("J122").Value = 2982.3333
$ws.Range("K122").Value = 4857.2856
$ws.Range("L122").Value = 8946.999899999999
$ws.Range("M122").Value = -2407.2856
$ws.Range("N122").Value = -13846.9999
$ws.Range("H126").Value = 1954.2609
$ws.Range("I126").Value = 1983.2858
$ws.Range("J126").Value = 1649.5
$ws.Range("K126").Value = 5949.857400000001
$ws.Range("L126").Value = 4948.5
$ws.Range("M126").Value = -3479.857400000001
$ws.Range("N126").Value = -9888.5
$ws.Range("H132").Value = 21121214
$ws.Range("I132").Value = 22294172
$ws.Range("K132").Value = 66882516
$ws.Range("M132").Value = -66879986
$ws.Range("H136").Value = 15884894
$ws.Range("I136").Value = 18530990
$ws.Range("J136").Value = 8323
$ws.Range("K136").Value = 55592970
$ws.Range("L136").Value = 24969
$ws.Range("M136").Value = -55590420
$ws.Range("N136").Value = -30069
